# Apply BackLog.xlsx update:
#  - Row 3 (Historia: Registrarme en la app) comments/requirements text updated
#    to reflect a registration form with 11 Edit Texts (adding usuario/contraseña
#    fields) instead of 8.
#  - Row 3 height grows to fit the longer wrapped text.
#  - The active window scroll position / selection is reset (was scrolled down
#    to H12, now back at top with G4 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Requerimientos" (G3) and "Comentarios" (I3) text for the
# --- "Registrarme en la app" story row ---
# (I3 set first so the shared-string table gains the two new entries in the
# same order the source workbook ended up with: "Se solicitara..." then
# "Activity con formulario...")
$ws.Range("I3").Value = "Se solicitara Nombre, Apellido, Sexo, correo electrónico, número de identificación, nombre del negocio, fecha de nacimiento, dirección del negocio, usuario y contraseña y repita constraseña"
$ws.Range("G3").Value = 'Activity con formulario de registro con 11 Edit Text con hint para ingresar los datos y boton de colo azúl para enviar informacion, 1 Toast que diga "usuario registrado satisfactoriamente, Bienvenido"'

# --- Grow row 3 so the longer wrapped text still fits ---
$ws.Rows(3).RowHeight = 129.6

# --- Reset view: scroll back to the top-left and select G4 (was topLeftCell
# --- A9 / selection H12) ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("G4").Select()
